$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Two new quiz entries (rows 19 & 20): Da Yu Ding (id 17) and
# Da Ke Ding (id 18).
# ------------------------------------------------------------------

# Numeric cells first (id / start_time / end_time) - these never touch
# the shared-string table.
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("G19").Value = -1040
$ws.Range("H19").Value = -1003
$ws.Range("G20").Value = -950
$ws.Range("H20").Value = -886

# Text cells, written in the exact order needed to reproduce the
# original shared-string insertion sequence.
$ws.Range("B19").Value = "大盂鼎"
$ws.Range("B20").Value = "大克鼎"
$ws.Range("C20").Value = "Da Ke Ding"
$ws.Range("J20").Value = "https://www.shanghaimuseum.net/resource/museum_files/show_files/20151104094055028/index.html"
$ws.Range("D20").Value = "Ritual Cauldron (Ding) of Duke Ke"
$ws.Range("E20").Value = "Fufeng, Shaanxi"
$ws.Range("I20").Value = "static/images/大克鼎.png"
$ws.Range("I19").Value = "static/images/大盂鼎.png"
$ws.Range("C19").Value = "Da Yu Ding"
$ws.Range("D19").Value = "Ritual Cauldron (Ding) of Duke Yu"
$ws.Range("J19").Value = "https://www.chnmuseum.cn/portals/0/web/zt/202106dayuding/"

# Cells that reuse strings already present in the workbook.
$ws.Range("E19").Value = "Qishan, Shaanxi"
$ws.Range("F19").Value = "National Museum of China"
$ws.Range("F20").Value = "Shanghai Museum"

# Make J20 a real hyperlink (like the other url cells that already
# carry one), then restore the shared "hyperlink" cell style used
# elsewhere in the sheet (e.g. J8) so J20 matches visually.
$ws.Hyperlinks.Add($ws.Range("J20"), "https://www.shanghaimuseum.net/resource/museum_files/show_files/20151104094055028/index.html")
$ws.Range("J8").Copy()
$ws.Range("J20").PasteSpecial(-4122)

# ------------------------------------------------------------------
# View tweaks: zoom level and active selection.
# ------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 81
$ws.Range("C15").Select()

$wb.Save()
